$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 4 data rows (rows 2-5). This shifts every subsequent
# row up by 4, dropping the final 4 rows of the series (old rows 39-42)
# and shrinking the used range from A1:E42 to A1:E38.
$ws.Rows("2:5").Delete()
